# Apply updated cryptocurrency price/volume figures (and the handful of
# rows whose coin/link pairs were re-sorted) from the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.183.89'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '''1.814.49'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("D4").Value = '''1.004'
$ws.Range("D5").Value = '''340.57'
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("D6").Value = '''1.005'
$ws.Range("E6").Value = '  +0.54%  '
$ws.Range("D7").Value = '''0.3919'
$ws.Range("E7").Value = '  +2.96%  '
$ws.Range("D8").Value = '''0.3489'
$ws.Range("E8").Value = '  +0.61%  '
$ws.Range("D9").Value = '''48.37'
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("D10").Value = '''1.196'
$ws.Range("E10").Value = '  -0.86%  '
$ws.Range("D11").Value = '''0.07573'
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").Value = '''1.009'
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").Value = '''22.13'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").Value = '''6.525'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '''1.818.35'
$ws.Range("E15").Value = '  +1.20%  '
$ws.Range("D16").Value = '''7.188'
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("D17").Value = '''0.00001106'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '''0.06739'
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("D19").Value = '''85.26'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").Value = '''1.003'
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").Value = '''17.81'
$ws.Range("E21").Value = '  +2.24%  '
$ws.Range("D22").Value = '''6.574'
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").Value = '''28.057.54'
$ws.Range("E23").Value = '  +2.26%  '
$ws.Range("D24").Value = '''12.51'
$ws.Range("E24").Value = '  -0.69%  '
$ws.Range("D25").Value = '''2.424'
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").Value = '''1.494'
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '''2.532'
$ws.Range("E27").Value = '  -1.84%  '
$ws.Range("D28").Value = '''21.31'
$ws.Range("E28").Value = '  -1.15%  '
$ws.Range("D29").Value = '''154.33'
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("D30").Value = '''2.033.37'
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("D31").Value = '''136.20'
$ws.Range("E31").Value = '  +1.44%  '
$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").Value = '''4.057'
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''6.160'
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("D34").Value = '''0.08764'
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("D35").Value = '''13.06'
$ws.Range("E35").Value = '  -2.16%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '''0.06574'
$ws.Range("E36").Value = '  +2.77%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '''5.479'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.02428'
$ws.Range("B39").Value = 'WEMIXTOKEN'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").Value = '''1.618'
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.6914'
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("D41").Value = '''0.2224'
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("D42").Value = '''1.263'
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").Value = '''8.539'
$ws.Range("E43").Value = '  -4.45%  '
$ws.Range("D44").Value = '''14.70'
$ws.Range("E44").Value = '  +1.88%  '
$ws.Range("D45").Value = '''0.6478'
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").Value = '''3.874'
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").Value = '''2.155'
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").Value = '''131.40'
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("D49").Value = '''0.07202'
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").Value = '''80.19'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = '''1.252'
$ws.Range("E51").Value = '  +2.33%  '
